# Adding 'Perfadex Plus' as a permissible value for perfusion and transport solution

$wb = $excel.ActiveWorkbook

# --- perfusion_solution sheet: insert 'Perfadex Plus' right after 'UWS' (new row 2) ---
$wsPerf = $wb.Worksheets.Item("perfusion_solution")
$wsPerf.Rows.Item(2).Insert()
$wsPerf.Cells.Item(2, 1).Value = "Perfadex Plus"
$wsPerf.Cells.Item(2, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000249"

# --- transport_solution sheet: insert 'Perfadex Plus' right after 'DMEM' (new row 5) ---
$wsTrans = $wb.Worksheets.Item("transport_solution")
$wsTrans.Rows.Item(5).Insert()
$wsTrans.Cells.Item(5, 1).Value = "Perfadex Plus"
$wsTrans.Cells.Item(5, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000249"

# --- Organ sheet: extend the data validation ranges to include the new row ---
$wsOrgan = $wb.Worksheets.Item("Organ")

$rngPerf = $wsOrgan.Range("E2:E1001")
$rngPerf.Validation.Modify(3, 1, 1, "='perfusion_solution'!`$A`$1:`$A`$7")

$rngTrans = $wsOrgan.Range("F2:F1001")
$rngTrans.Validation.Modify(3, 1, 1, "='transport_solution'!`$A`$1:`$A`$12")

# --- .metadata sheet: bump pav:createdOn timestamp ---
$wsMeta = $wb.Worksheets.Item(".metadata")
$wsMeta.Cells.Item(2, 3).Value = "2023-09-01T13:52:59-07:00"
